$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "275 TOPS (INT8-Sparse)"
$ws.Range("C3").Value = "275 TOPS (INT8-Sparse)"
$ws.Range("C4").Value = "248 TOPS (INT8-Sparse)"
$ws.Range("C5").Value = "200 TOPS (INT8-Sparse)"
$ws.Range("C6").Value = "157 TOPS (INT8-Sparse)"
$ws.Range("C7").Value = "117 TOPS (INT8-Sparse)"
$ws.Range("C8").Value = "67 TOPS (INT8-Sparse)"
$ws.Range("C9").Value = "67 TOPS (INT8-Sparse)"
$ws.Range("C10").Value = "34 TOPS (INT8-Sparse)"
$ws.Range("C14").Value = "30 TOPS (INT8-Sparse)"
$ws.Range("C15").Value = "32 TOPS (INT8-Sparse)"
$ws.Range("C16").Value = "32 TOPS (INT8-Sparse)"
$ws.Range("C17").Value = "21 TOPS (INT8-Sparse)"
$ws.Range("C18").Value = "21 TOPS (INT8-Sparse)"
$ws.Range("C19").Value = "1.26 TFLOPS (FP16-Dense)"
$ws.Range("C20").Value = "1.33 TFLOPS (FP16-Dense)"
$ws.Range("C21").Value = "1.33 TFLOPS (FP16-Dense)"
$ws.Range("C22").Value = "1.33 TFLOPS (FP16-Dense)"
$ws.Range("C23").Value = "0.472 TFLOPS (FP16-Dense)"
